$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B54").Value = "Cụm 7"
$ws.Range("B50").Value = "Tổ dân phố 1 (Phường Sơn Lộc)"
